# ---------------------------------------------------------------------------
# Adds the "Prob(success)" efficiency-index column (I) plus a couple of
# "divided by N" notes (J) to the "R input" sheet, per the commit:
#   "getting efficiency index correct. finalized on efficiency =
#    Time/Strikes x probability(success). Probability of success = number
#    of successful splits out of all cobble experiments."
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R input")
$ws.Activate()

# --- Cosmetic window placement (best effort; may not round-trip) ----------
try { $wb.Windows.Item(1).Left = 1920 } catch {}

# --- Header row -------------------------------------------------------------
$ws.Range("I1").Value = "Prob(success)"
$ws.Range("I1").Font.Bold = $true

$ws.Range("J1").Value = "Notes"
$ws.Range("J1").Font.Bold = $true

# --- Group E3 (rows 2-34): Prob(success) = 14/33 ---------------------------
$ws.Range("I2").Formula = "=14/33"
$ws.Range("J2").Value = "divided by 33"
$ws.Range("I3").Formula = "=14/33"
$ws.Range("I4:I34").Formula = "=14/33"

# --- Group E2 (rows 35-70): Prob(success) = 22/36 --------------------------
$ws.Range("I35").Formula = "=22/36"
$ws.Range("J35").Value = "divided by 36"
$ws.Range("I36").Formula = "=22/36"
$ws.Range("I37:I69").Formula = "=22/36"
$ws.Range("I70").Formula = "=22/36"

# --- Group E1 (rows 71-103): Prob(success) = 10/33 --------------------------
$ws.Range("I71").Formula = "=10/33"
$ws.Range("J71").Value = "divided by 33"
$ws.Range("I72").Formula = "=10/33"
$ws.Range("I73:I103").Formula = "=10/33"

# --- Groups N10+N7 (rows 104-123): Prob(success) = 2/10 ---------------------
$ws.Range("I104").Formula = "=2/10"
$ws.Range("I105").Formula = "=2/10"
$ws.Range("I106:I123").Formula = "=2/10"

# --- Group N1 (rows 124-133): Prob(success) = 0.4 (typed literal) ----------
$ws.Range("I124:I133").Value = 0.4

# --- Group N6 (rows 134-143): Prob(success) = 1/10, interspersed literals --
$ws.Range("I134").Formula = "=1/10"
$ws.Range("I135:I136").Value = 0.1
$ws.Range("I137:I143").Formula = "=1/10"
$ws.Range("I138:I139").Value = 0.1
$ws.Range("I140:I143").Formula = "=1/10"
$ws.Range("I141:I142").Value = 0.1
$ws.Range("I143").Formula = "=1/10"

# --- Group N4 (rows 144-153): Prob(success) = 0.3 (typed literal) ----------
$ws.Range("I144:I153").Value = 0.3

# --- Group N9 (rows 154-163): Prob(success) = 0.6 (typed literal) ----------
$ws.Range("I154:I163").Value = 0.6

# --- Group N5 (rows 164-173): Prob(success) = 0 (typed literal) ------------
$ws.Range("I164:I173").Value = 0

# --- Group N11 (rows 174-183): Prob(success) = 0.2 (typed literal) ---------
$ws.Range("I174:I183").Value = 0.2

# --- Group N2 (rows 184-193): Prob(success) = 0 (typed literal) ------------
$ws.Range("I184:I193").Value = 0

# --- View state: scrolled near the bottom, last cell selected --------------
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 165
    $win.ScrollColumn = 1
} catch {}
$ws.Range("K192").Select()

Write-Output "Added Prob(success) column I1:I193 and notes J1/J2/J35/J71"
